$d = $word.ActiveDocument

function Append-ToParagraphEnd($paraIndex, $text) {
    $para = $d.Paragraphs($paraIndex)
    $r = $para.Range
    $insertPos = $r.End - 1
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter($text)
}

# ------------------------------------------------------------------
# 1) "Needs:" paragraph -- append " accuracy, efficiency." after the
#    existing "...reliability," sentence.
# ------------------------------------------------------------------
Append-ToParagraphEnd 3 " accuracy, efficiency."

# ------------------------------------------------------------------
# 2) "Different Chips Alternatives" > "Microcontrollers:" section
#    Pros: line -> append "Easy to use, lower power consumption"
#    Cons: line -> append " Serial execution, internal peripherals can limit scope."
# ------------------------------------------------------------------
Append-ToParagraphEnd 25 "Easy to use, lower power consumption"
Append-ToParagraphEnd 26 " Serial execution, internal peripherals can limit scope."

# ------------------------------------------------------------------
# 3) "FPGA:" section
#    Pros: line -> "Pros: Read in data in parrell" becomes
#                  "Pros: Can be programmed at logic level (parallel processing)."
#    Cons: line -> append " More complex coding, no control of power consumption."
# ------------------------------------------------------------------
$fpgaProsPara = $d.Paragraphs(29)
$fpgaProsRange = $fpgaProsPara.Range
$fpgaProsRange.Find.Execute("Read in data in parrell", $false, $false, $false, $false, $false, $true, 1, $false, "Can be programmed at logic level (parallel processing).", 2)

Append-ToParagraphEnd 30 " More complex coding, no control of power consumption."
